$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column from 2023-09-15 (45184) to 2023-09-16 (45185)
# for rows 2 through 5 (C2:C5), preserving the existing date number format.
$newDateSerial = 45185

$ws.Range("C2").Value = $newDateSerial
$ws.Range("C3").Value = $newDateSerial
$ws.Range("C4").Value = $newDateSerial
$ws.Range("C5").Value = $newDateSerial
